# Applies the diff between the original and the updated
# StructureDefinition-biospecimen-laterality.xlsx workbook.
#
# Net visible changes (the shared-string re-indexing/dedup seen in the raw
# OOXML diff is just a side effect of these content edits):
#   Metadata sheet:
#     B8  (Date)          2025-05-21T14:22:51+00:00  -> 2025-06-13T15:45:04+00:00
#     B15 (FHIR Version)  4.3.0                       -> 4.0.1
#   Elements sheet:
#     K3  (Extension.id Type(s))           id\n        -> string\n
#     AJ2 (Extension Constraint(s))        drop the "unless an empty Parameters
#                                           resource {...} or $this is Parameters"
#                                           clause from the ele-1 constraint
#     M6  (Extension.value[x] Definition)  R4B -> R4 in the Extensibility link

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$meta.Range("B15").Value = "4.0.1"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K3").Value = "string`n"
$elements.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
$elements.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
